# Apply updated crypto market data to the worksheet, matching the upstream
# GitHub Actions data-refresh commit. Column D holds price text and column E
# holds 1h volume-change text; some prices are plain numeric-looking strings
# (e.g. "317.98") that Excel would otherwise auto-convert to numbers, so we
# pre-format those specific cells as Text and restore the default style after
# writing, keeping the cells' appearance identical to the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value would otherwise be parsed as a number by Excel;
# mark them as Text first so they are stored as strings, like the source data.
$textCells = @(
    "D5", "D6", "D7", "D10", "D13", "D14", "D17", "D19",
    "D20", "D22", "D23", "D24", "D25", "D27", "D30", "D31",
    "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40",
    "D41", "D43", "D44", "D50", "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '45.221.75'
$ws.Cells.Item(2, 5).Value = '  +4.88%  '
$ws.Cells.Item(3, 4).Value = '2.432.62'
$ws.Cells.Item(3, 5).Value = '  +1.85%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = '317.98'
$ws.Cells.Item(5, 5).Value = '  +4.63%  '
$ws.Cells.Item(6, 4).Value = '103.44'
$ws.Cells.Item(6, 5).Value = '  +7.07%  '
$ws.Cells.Item(7, 4).Value = '0.515'
$ws.Cells.Item(7, 5).Value = '  +2.06%  '
$ws.Cells.Item(8, 5).Value = '  -0.08%  '
$ws.Cells.Item(9, 5).Value = '  +9.03%  '
$ws.Cells.Item(10, 4).Value = '35.76'
$ws.Cells.Item(10, 5).Value = '  +2.57%  '
$ws.Cells.Item(11, 5).Value = '  +1.43%  '
$ws.Cells.Item(12, 5).Value = '  -2.17%  '
$ws.Cells.Item(13, 4).Value = '18.39'
$ws.Cells.Item(13, 5).Value = '  -0.71%  '
$ws.Cells.Item(14, 4).Value = '6.98'
$ws.Cells.Item(14, 5).Value = '  +2.28%  '
$ws.Cells.Item(15, 4).Value = '2.814.74'
$ws.Cells.Item(15, 5).Value = '  +2.48%  '
$ws.Cells.Item(16, 4).Value = '2.437.49'
$ws.Cells.Item(16, 5).Value = '  +2.65%  '
$ws.Cells.Item(17, 4).Value = '0.831'
$ws.Cells.Item(17, 5).Value = '  +2.34%  '
$ws.Cells.Item(18, 4).Value = '45.099.69'
$ws.Cells.Item(18, 5).Value = '  +4.64%  '
$ws.Cells.Item(19, 4).Value = '12.26'
$ws.Cells.Item(19, 5).Value = '  +2.14%  '
$ws.Cells.Item(20, 4).Value = '6.37'
$ws.Cells.Item(20, 5).Value = '  +0.25%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0925'
$ws.Cells.Item(21, 5).Value = '  +3.82%  '
$ws.Cells.Item(22, 4).Value = '69.00'
$ws.Cells.Item(22, 5).Value = '  +1.20%  '
$ws.Cells.Item(23, 4).Value = '243.44'
$ws.Cells.Item(23, 5).Value = '  +2.76%  '
$ws.Cells.Item(24, 4).Value = '2.28'
$ws.Cells.Item(24, 5).Value = '  +1.52%  '
$ws.Cells.Item(25, 4).Value = '2.51'
$ws.Cells.Item(25, 5).Value = '  +2.73%  '
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).Value = '25.37'
$ws.Cells.Item(27, 5).Value = '  +3.03%  '
$ws.Cells.Item(28, 5).Value = '  -7.97%  '
$ws.Cells.Item(29, 5).Value = '  +1.82%  '
$ws.Cells.Item(30, 2).Value = 'OKB'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(30, 4).Value = '49.32'
$ws.Cells.Item(30, 5).Value = '  +2.63%  '
$ws.Cells.Item(31, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(31, 4).Value = '33.05'
$ws.Cells.Item(31, 5).Value = '  +3.09%  '
$ws.Cells.Item(32, 4).Value = '20.33'
$ws.Cells.Item(32, 5).Value = '  +12.89%  '
$ws.Cells.Item(33, 4).Value = '0.127'
$ws.Cells.Item(33, 5).Value = '  +10.43%  '
$ws.Cells.Item(34, 4).Value = '5.23'
$ws.Cells.Item(34, 5).Value = '  +2.68%  '
$ws.Cells.Item(35, 5).Value = '  +0.22%  '
$ws.Cells.Item(36, 4).Value = '0.0766'
$ws.Cells.Item(36, 5).Value = '  +3.54%  '
$ws.Cells.Item(37, 4).Value = '1.89'
$ws.Cells.Item(37, 5).Value = '  +2.74%  '
$ws.Cells.Item(38, 4).Value = '4.48'
$ws.Cells.Item(38, 5).Value = '  +3.46%  '
$ws.Cells.Item(39, 4).Value = '2.85'
$ws.Cells.Item(39, 5).Value = '  -0.25%  '
$ws.Cells.Item(40, 4).Value = '124.23'
$ws.Cells.Item(40, 5).Value = '  -3.65%  '
$ws.Cells.Item(41, 4).Value = '0.109'
$ws.Cells.Item(41, 5).Value = '  +1.85%  '
$ws.Cells.Item(42, 5).Value = '  -2.46%  '
$ws.Cells.Item(43, 4).Value = '21.13'
$ws.Cells.Item(43, 5).Value = '  -0.45%  '
$ws.Cells.Item(44, 4).Value = '0.0290'
$ws.Cells.Item(44, 5).Value = '  +3.79%  '
$ws.Cells.Item(45, 4).Value = '1.937.95'
$ws.Cells.Item(45, 5).Value = '  +0.17%  '
$ws.Cells.Item(46, 5).Value = '  +4.98%  '
$ws.Cells.Item(47, 5).Value = '  -2.01%  '
$ws.Cells.Item(48, 5).Value = '  +0.00%  '
$ws.Cells.Item(49, 5).Value = '  +16.24%  '
$ws.Cells.Item(50, 4).Value = '76.28'
$ws.Cells.Item(50, 5).Value = '  +6.24%  '
$ws.Cells.Item(51, 4).Value = '4.73'
$ws.Cells.Item(51, 5).Value = '  +6.71%  '

# Restore the default (style-less) formatting on the text-forced cells so the
# only observable change is the cell content, matching the source diff.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
